$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Entity types table")
Write-Host $ws1.Name
